# Macroferia Regional de Talca - Zapallo: weekly fruit/hortaliza update.
# Inserts two new weekly records at rows 260-261 (pushing the existing
# rows 260-284 down to 262-286), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 260; everything at/after
# row 260 (through 284) shifts down to 262..286.
$ws.Rows("260:261").Insert()

# ---- New row 260 ----
$ws.Range("A260").Value = 5
$ws.Range("B260").Value = "Macroferia Regional de Talca"
$ws.Range("C260").Value = "Maule"
$ws.Range("D260").Value = 44783
$ws.Range("E260").Value = 7
$ws.Range("F260").Value = 100112045
$ws.Range("G260").Value = "Zapallo"
$ws.Range("H260").Value = "Camote"
$ws.Range("I260").Value = "1a (guarda)"
$ws.Range("J260").Value = 900
$ws.Range("K260").Value = 650
$ws.Range("L260").Value = 650
$ws.Range("M260").Value = 650
$ws.Range("N260").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O260").Value = "Región del Maule"
$ws.Range("P260").Value = 650
$ws.Range("Q260").Value = 1
$ws.Range("R260").Value = "Hortaliza"

# ---- New row 261 ----
$ws.Range("A261").Value = 5
$ws.Range("B261").Value = "Macroferia Regional de Talca"
$ws.Range("C261").Value = "Maule"
$ws.Range("D261").Value = 44783
$ws.Range("E261").Value = 7
$ws.Range("F261").Value = 100112045
$ws.Range("G261").Value = "Zapallo"
$ws.Range("H261").Value = "Paine"
$ws.Range("I261").Value = "1a (guarda)"
$ws.Range("J261").Value = 2000
$ws.Range("K261").Value = 230
$ws.Range("L261").Value = 230
$ws.Range("M261").Value = 230
$ws.Range("N261").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O261").Value = "Región del Maule"
$ws.Range("P261").Value = 230
$ws.Range("Q261").Value = 1
$ws.Range("R261").Value = "Hortaliza"
